$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values, updated per the automatic electricity price refresh
$ws.Range("A2").Value = 45944
$ws.Range("B2").Value = 116.29
$ws.Range("C2").Value = 109.22
$ws.Range("D2").Value = 107.06
$ws.Range("E2").Value = 106.93
$ws.Range("F2").Value = 107.22
$ws.Range("G2").Value = 108.5
$ws.Range("H2").Value = 115.58
$ws.Range("I2").Value = 126.63
$ws.Range("J2").Value = 155.86
$ws.Range("K2").Value = 122.51
$ws.Range("L2").Value = 105.43
$ws.Range("M2").Value = 88.31999999999999
$ws.Range("N2").Value = 81.83
$ws.Range("O2").Value = 75.13
$ws.Range("P2").Value = 74.25
$ws.Range("Q2").Value = 77.63
$ws.Range("R2").Value = 81.93000000000001
$ws.Range("S2").Value = 92.55
$ws.Range("T2").Value = 112.33
$ws.Range("U2").Value = 145.16
$ws.Range("V2").Value = 160.2
$ws.Range("W2").Value = 144.19
$ws.Range("X2").Value = 130.07
$ws.Range("Y2").Value = 111.39
$ws.Range("Z2").Value = 110.68

$ws.Range("AB2").Value = 136.46
$ws.Range("AD2").Value = 152.2
$ws.Range("AF2").Value = 139.18
$ws.Range("AG2").Value = "1h-17h"
